# Auto-generated Excel COM-interop script to apply the scheduled runner update
# to the Seraph Profits market-data workbook across all 8 job sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 814.1429000000001
$ws.Range("I32").Value2 = 814.1429000000001
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 814.1429000000001
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -488.1429000000001
$ws.Range("N32").ClearContents()
$ws.Range("H86").Value2 = 1500
$ws.Range("I86").Value2 = 1500
$ws.Range("K86").Value2 = 1500
$ws.Range("M86").Value2 = -377
$ws.Range("H89").Value2 = 1500
$ws.Range("I89").Value2 = 1500
$ws.Range("K89").Value2 = 7500
$ws.Range("M89").Value2 = -1884
$ws.Range("H138").Value2 = 3040.5715
$ws.Range("J138").Value2 = 8231.166999999999
$ws.Range("L138").Value2 = 24693.501
$ws.Range("N138").Value2 = -34973.501

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 23810774
$ws.Range("I2").Value2 = 27778420
$ws.Range("J2").Value2 = 4900
$ws.Range("K2").Value2 = 27778420
$ws.Range("L2").Value2 = 4900
$ws.Range("M2").Value2 = -27778307
$ws.Range("N2").Value2 = -5126
$ws.Range("H116").Value2 = 23810774
$ws.Range("I116").Value2 = 27778420
$ws.Range("J116").Value2 = 4900
$ws.Range("K116").Value2 = 27778420
$ws.Range("L116").Value2 = 4900
$ws.Range("M116").Value2 = -27776126
$ws.Range("N116").Value2 = -9488
$ws.Range("H122").Value2 = 1260975.6
$ws.Range("I122").Value2 = 1260975.6
$ws.Range("K122").Value2 = 3782926.8
$ws.Range("M122").Value2 = -3780476.8

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 23810774
$ws.Range("I3").Value2 = 27778420
$ws.Range("J3").Value2 = 4900
$ws.Range("K3").Value2 = 27778420
$ws.Range("L3").Value2 = 4900
$ws.Range("M3").Value2 = -27778306
$ws.Range("N3").Value2 = -5128
$ws.Range("H105").Value2 = 4905717
$ws.Range("I105").Value2 = 10420663
$ws.Range("K105").Value2 = 10420663
$ws.Range("M105").Value2 = -10418916
$ws.Range("H107").Value2 = 1796.8572
$ws.Range("I107").Value2 = 1742.7693
$ws.Range("K107").Value2 = 1742.7693
$ws.Range("M107").Value2 = 177.2307000000001
$ws.Range("H123").Value2 = 0
$ws.Range("J123").Value2 = 0
$ws.Range("L123").Value2 = 0
$ws.Range("N123").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value2 = 20009
$ws.Range("I15").Value2 = 0
$ws.Range("J15").Value2 = 20009
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 20009
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value2 = -20349
$ws.Range("H29").Value2 = 7118.5
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 7118.5
$ws.Range("K29").Value2 = 0
$ws.Range("L29").Value2 = 7118.5
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value2 = -7704.5
$ws.Range("H99").Value2 = 10061.871
$ws.Range("I99").Value2 = 5809.0557
$ws.Range("K99").Value2 = 5809.0557
$ws.Range("M99").Value2 = -4311.0557
$ws.Range("H107").Value2 = 83333896
$ws.Range("I107").Value2 = 125000270
$ws.Range("J107").Value2 = 1149
$ws.Range("K107").Value2 = 125000270
$ws.Range("L107").Value2 = 1149
$ws.Range("M107").Value2 = -124998350
$ws.Range("N107").Value2 = -4989
$ws.Range("H122").Value2 = 1260
$ws.Range("J122").Value2 = 1400
$ws.Range("L122").Value2 = 4200
$ws.Range("N122").Value2 = -9100
$ws.Range("H126").Value2 = 10061.871
$ws.Range("I126").Value2 = 5809.0557
$ws.Range("K126").Value2 = 17427.1671
$ws.Range("M126").Value2 = -14957.1671
$ws.Range("H132").Value2 = 2214.5881
$ws.Range("I132").Value2 = 1858.1538
$ws.Range("K132").Value2 = 5574.4614
$ws.Range("M132").Value2 = -3044.4614

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 1731016.1
$ws.Range("J4").Value2 = 934.5714
$ws.Range("L4").Value2 = 2803.7142
$ws.Range("N4").Value2 = -3027.7142
$ws.Range("H121").Value2 = 714.8333
$ws.Range("I121").Value2 = 773
$ws.Range("J121").Value2 = 598.5
$ws.Range("K121").Value2 = 2319
$ws.Range("L121").Value2 = 1795.5
$ws.Range("M121").Value2 = -1009
$ws.Range("N121").Value2 = -4415.5
$ws.Range("H131").Value2 = 981
$ws.Range("J131").Value2 = 980
$ws.Range("L131").Value2 = 2940
$ws.Range("N131").Value2 = -13020

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("N34").ClearContents()
$ws.Range("H70").Value2 = 4624.25
$ws.Range("I70").Value2 = 1250
$ws.Range("J70").Value2 = 7998.5
$ws.Range("K70").Value2 = 1250
$ws.Range("L70").Value2 = 7998.5
$ws.Range("M70").Value2 = -980
$ws.Range("N70").Value2 = -8538.5
$ws.Range("H73").Value2 = 4624.25
$ws.Range("I73").Value2 = 1250
$ws.Range("J73").Value2 = 7998.5
$ws.Range("K73").Value2 = 1250
$ws.Range("L73").Value2 = 7998.5
$ws.Range("M73").Value2 = -314
$ws.Range("N73").Value2 = -9870.5
$ws.Range("H76").Value2 = 0
$ws.Range("J76").Value2 = 0
$ws.Range("L76").Value2 = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value2 = 0
$ws.Range("J79").Value2 = 0
$ws.Range("L79").Value2 = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value2 = 4298.5386
$ws.Range("J80").Value2 = 5499.6665
$ws.Range("L80").Value2 = 5499.6665
$ws.Range("N80").Value2 = -7495.6665
$ws.Range("H83").Value2 = 4298.5386
$ws.Range("J83").Value2 = 5499.6665
$ws.Range("L83").Value2 = 27498.3325
$ws.Range("N83").Value2 = -37482.3325
$ws.Range("H123").Value2 = 34979.89
$ws.Range("J123").Value2 = 34979.89
$ws.Range("L123").Value2 = 34979.89
$ws.Range("N123").Value2 = -39879.89

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 1713.1428
$ws.Range("J7").Value2 = 1997.6666
$ws.Range("L7").Value2 = 1997.6666
$ws.Range("N7").Value2 = -2221.6666
$ws.Range("H34").Value2 = 15332
$ws.Range("J34").Value2 = 8000
$ws.Range("L34").Value2 = 8000
$ws.Range("N34").Value2 = -8344
$ws.Range("H43").Value2 = 364793.56
$ws.Range("I43").Value2 = 6505.5
$ws.Range("J43").Value2 = 424508.25
$ws.Range("K43").Value2 = 6505.5
$ws.Range("L43").Value2 = 424508.25
$ws.Range("M43").Value2 = -6312.5
$ws.Range("N43").Value2 = -424894.25
$ws.Range("H68").Value2 = 2850.25
$ws.Range("I68").Value2 = 2132.6667
$ws.Range("K68").Value2 = 2132.6667
$ws.Range("M68").Value2 = -1383.6667
$ws.Range("H71").Value2 = 2850.25
$ws.Range("I71").Value2 = 2132.6667
$ws.Range("K71").Value2 = 10663.3335
$ws.Range("M71").Value2 = -6919.333500000001
$ws.Range("H122").Value2 = 2076
$ws.Range("I122").Value2 = 1152.5
$ws.Range("K122").Value2 = 3457.5
$ws.Range("M122").Value2 = -1007.5
$ws.Range("H126").Value2 = 1713.1428
$ws.Range("J126").Value2 = 1997.6666
$ws.Range("L126").Value2 = 5992.9998
$ws.Range("N126").Value2 = -10932.9998

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 5826.722
$ws.Range("I62").Value2 = 3677.8
$ws.Range("J62").Value2 = 6653.231
$ws.Range("K62").Value2 = 3677.8
$ws.Range("L62").Value2 = 6653.231
$ws.Range("M62").Value2 = -3053.8
$ws.Range("N62").Value2 = -7901.231
$ws.Range("H64").Value2 = 56797
$ws.Range("J64").Value2 = 56797
$ws.Range("L64").Value2 = 56797
$ws.Range("N64").Value2 = -57293
$ws.Range("H65").Value2 = 5826.722
$ws.Range("I65").Value2 = 3677.8
$ws.Range("J65").Value2 = 6653.231
$ws.Range("K65").Value2 = 18389
$ws.Range("L65").Value2 = 33266.155
$ws.Range("M65").Value2 = -15269
$ws.Range("N65").Value2 = -39506.155
$ws.Range("H67").Value2 = 56797
$ws.Range("J67").Value2 = 56797
$ws.Range("L67").Value2 = 56797
$ws.Range("N67").Value2 = -58513
$ws.Range("H81").Value2 = 2552.2307
$ws.Range("I81").Value2 = 2707.182
$ws.Range("J81").Value2 = 1700
$ws.Range("K81").Value2 = 5414.364
$ws.Range("L81").Value2 = 3400
$ws.Range("M81").Value2 = -4353.364
$ws.Range("N81").Value2 = -5522
$ws.Range("H84").Value2 = 2552.2307
$ws.Range("I84").Value2 = 2707.182
$ws.Range("J84").Value2 = 1700
$ws.Range("K84").Value2 = 27071.82
$ws.Range("L84").Value2 = 17000
$ws.Range("M84").Value2 = -21767.82
$ws.Range("N84").Value2 = -27608
$ws.Range("H94").Value2 = 37992.4
$ws.Range("I94").Value2 = 37992.4
$ws.Range("J94").Value2 = 0
$ws.Range("K94").Value2 = 37992.4
$ws.Range("L94").Value2 = 0
$ws.Range("M94").Value2 = -37091.4
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value2 = 2616.5
$ws.Range("I100").Value2 = 2418.5715
$ws.Range("K100").Value2 = 4837.143
$ws.Range("M100").Value2 = -4296.143
